$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.862.22'
$ws.Range('E2').Value = '  -1.36%  '

$ws.Range('D3').Value = '1.755.14'
$ws.Range('E3').Value = '  -4.05%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.007'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.89%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '337.67'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.66%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.36%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3754'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -4.81%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3342'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -4.62%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '45.86'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -4.79%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.115'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -7.43%  '

$ws.Range('E11').Value = '  -6.58%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.005'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.99%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.03'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.05%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.148'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -6.39%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.127'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.57%  '

$ws.Range('D16').Value = '1.757.27'
$ws.Range('E16').Value = '  -3.81%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001045'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -6.04%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06556'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.26%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '79.93'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -6.85%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.54%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.81'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -6.69%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.240'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -5.58%  '

$ws.Range('D23').Value = '27.914.54'
$ws.Range('E23').Value = '  -1.10%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.60'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -9.43%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.399'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.38%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '151.98'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.10%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.61'
$ws.Range('D27').Style = "Normal"

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.301'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -11.69%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.284'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -15.43%  '

$ws.Range('D30').Value = '1.960.33'
$ws.Range('E30').Value = '  -3.58%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '131.27'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.61%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.020'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.34%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.730'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -8.71%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08721'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.02%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '12.10'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -9.41%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02321'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -5.48%  '

$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.6515'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -7.15%  '

$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06155'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -6.60%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.100'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -8.55%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2097'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.06%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.208'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.90%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.447'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -10.22%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.980'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -7.76%  '

$ws.Range('E44').Value = '  +0.60%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.69'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -7.29%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.828'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.65%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5983'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -8.53%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '128.62'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.89%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.993'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -8.69%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07186'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.51%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.170'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.65%  '
